$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: new column F ---
$ws.Range("F1").Value = "Trening"
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats: copy A1's header style (bold, centered, bordered)

# --- Data rows 2-13: Timestamp (A), Seconds (B), Velocity (C), Acceleration_SMA (D), Velocity_Bin (E), Trening (F) ---
$rows = @(
    @{ r = 2;  a = 45685.64700069444; b = 1015.8; c = 13.8;   d = 3.401492255074637;  e = "10-15"; f = "Duża Gra" },
    @{ r = 3;  a = 45685.65112916667; b = 1372.5; c = 13.85;  d = 3.645100150789532;  e = "10-15"; f = "Duża Gra" },
    @{ r = 4;  a = 45685.65477962963; b = 1687.9; c = 12.84;  d = 3.595331192016603;  e = "10-15"; f = "Duża Gra" },
    @{ r = 5;  a = 45685.64699722223; b = 1015.5; c = 9.94;   d = 3.164574929646085;  e = "5-10";  f = "Duża Gra" },
    @{ r = 6;  a = 45685.64967314815; b = 1246.7; c = 9.68;   d = 3.149442638669694;  e = "5-10";  f = "Duża Gra" },
    @{ r = 7;  a = 45685.65112569444; b = 1372.2; c = 9.52;   d = 3.08819692475455;   e = "5-10";  f = "Duża Gra" },
    @{ r = 8;  a = 45685.67058518519; b = 3053.5; c = 11.7;   d = 2.875223670686994;  e = "10-15"; f = "Mała Gra" },
    @{ r = 9;  a = 45685.67079583334; b = 3071.7; c = 10.96;  d = 2.901540960584366;  e = "10-15"; f = "Mała Gra" },
    @{ r = 10; a = 45685.67585023148; b = 3508.4; c = 10.13;  d = 2.468041368893214;  e = "10-15"; f = "Mała Gra" },
    @{ r = 11; a = 45685.67079467593; b = 3071.6; c = 9.84;   d = 2.850619895117621;  e = "5-10";  f = "Mała Gra" },
    @{ r = 12; a = 45685.67584907408; b = 3508.3; c = 9.33;   d = 2.514835017068044;  e = "5-10";  f = "Mała Gra" },
    @{ r = 13; a = 45685.67589421296; b = 3512.2; c = 9.09;   d = 2.375224385942732;  e = "5-10";  f = "Mała Gra" }
)

# Register the two custom date numFmts (164 lowercase, 165 uppercase) on the
# first data cell, then just reapply the uppercase format (165) to the
# remaining date cells so no stray/unused style entries pile up.
$ws.Cells.Item(2, 1).Value = $rows[0].a
$ws.Cells.Item(2, 1).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(2, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

foreach ($row in $rows) {
    $r = $row.r
    if ($r -ne 2) {
        $ws.Cells.Item($r, 1).Value = $row.a
        $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    }
    $ws.Cells.Item($r, 2).Value = $row.b
    $ws.Cells.Item($r, 3).Value = $row.c
    $ws.Cells.Item($r, 4).Value = $row.d
    $ws.Cells.Item($r, 5).Value = $row.e
    $ws.Cells.Item($r, 6).Value = $row.f
}
